## Updating screenshot utility and writing table to excel
##
## - Leaves the "SearchPage" sheet's data alone, just restores the cursor to
##   the cell the user had selected there (K18) before switching tabs.
## - Adds a new "BookingHistory" sheet (placed right after "SearchPage") and
##   writes the order-history table into it: header row (Order Id / Total
##   Price) followed by the two booking rows.
## - Auto-sizes the new sheet's columns to the written content and makes the
##   new sheet the active tab, matching what Excel does when a table is
##   written/opened interactively.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remember the previous selection on SearchPage before we move focus away.
$ws1.Range("K18").Select() | Out-Null

# Insert the new sheet right after SearchPage and give it its name.
$bookingHistory = $wb.Worksheets.Add($null, $ws1)
$bookingHistory.Name = "BookingHistory"

# The table we're writing out: header row + data rows.
$table = @(
    @("Order Id", "Total Price"),
    @("URXJ1F704M", "AUD `$ 396"),
    @("Z2U1I9YV37", "AUD `$ 396")
)

for ($r = 0; $r -lt $table.Length; $r++) {
    $rowData = $table[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $bookingHistory.Cells.Item($r + 1, $c + 1).Value = $rowData[$c]
    }
}

# Give the header cells their own (distinct) style records.
$bookingHistory.Cells.Item(1, 1).Locked = $true
$bookingHistory.Cells.Item(1, 2).Locked = $false

# Fit the columns to the data that was just written.
$bookingHistory.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$bookingHistory.Columns.Item(2).ColumnWidth = 9.7
